# Auto update Excel log
# Appends newly-logged sensor rows (2026-01-30 PM) to the per-sensor
# history sheets, mirroring how the logger keeps extending each sheet.
$wb = $excel.ActiveWorkbook

# --- PIR: append 15 rows starting at row 60 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(60,1).NumberFormat = "@"
$ws.Cells.Item(60,1).Value = '2026-01-30'
$ws.Cells.Item(60,2).Value = '12:56:35'
$ws.Cells.Item(60,3).Value = '12:00'
$ws.Cells.Item(60,4).Value = 'Bathroom'
$ws.Cells.Item(60,5).Value = 'No Motion'
$ws.Cells.Item(60,6).Value = 'Inactive'
$ws.Cells.Item(61,1).NumberFormat = "@"
$ws.Cells.Item(61,1).Value = '2026-01-30'
$ws.Cells.Item(61,2).Value = '12:56:37'
$ws.Cells.Item(61,3).Value = '12:00'
$ws.Cells.Item(61,4).Value = 'Bathroom'
$ws.Cells.Item(61,5).Value = 'No Motion'
$ws.Cells.Item(61,6).Value = 'Inactive'
$ws.Cells.Item(62,1).NumberFormat = "@"
$ws.Cells.Item(62,1).Value = '2026-01-30'
$ws.Cells.Item(62,2).Value = '12:56:40'
$ws.Cells.Item(62,3).Value = '12:00'
$ws.Cells.Item(62,4).Value = 'Bathroom'
$ws.Cells.Item(62,5).Value = 'No Motion'
$ws.Cells.Item(62,6).Value = 'Inactive'
$ws.Cells.Item(63,1).NumberFormat = "@"
$ws.Cells.Item(63,1).Value = '2026-01-30'
$ws.Cells.Item(63,2).Value = '13:02:02'
$ws.Cells.Item(63,3).Value = '13:00'
$ws.Cells.Item(63,4).Value = 'Bathroom'
$ws.Cells.Item(63,5).Value = 'No Motion'
$ws.Cells.Item(63,6).Value = 'Inactive'
$ws.Cells.Item(64,1).NumberFormat = "@"
$ws.Cells.Item(64,1).Value = '2026-01-30'
$ws.Cells.Item(64,2).Value = '13:02:07'
$ws.Cells.Item(64,3).Value = '13:00'
$ws.Cells.Item(64,4).Value = 'Bathroom'
$ws.Cells.Item(64,5).Value = 'No Motion'
$ws.Cells.Item(64,6).Value = 'Inactive'
$ws.Cells.Item(65,1).NumberFormat = "@"
$ws.Cells.Item(65,1).Value = '2026-01-30'
$ws.Cells.Item(65,2).Value = '13:02:12'
$ws.Cells.Item(65,3).Value = '13:00'
$ws.Cells.Item(65,4).Value = 'Bathroom'
$ws.Cells.Item(65,5).Value = 'No Motion'
$ws.Cells.Item(65,6).Value = 'Inactive'
$ws.Cells.Item(66,1).NumberFormat = "@"
$ws.Cells.Item(66,1).Value = '2026-01-30'
$ws.Cells.Item(66,2).Value = '13:02:17'
$ws.Cells.Item(66,3).Value = '13:00'
$ws.Cells.Item(66,4).Value = 'Bathroom'
$ws.Cells.Item(66,5).Value = 'No Motion'
$ws.Cells.Item(66,6).Value = 'Inactive'
$ws.Cells.Item(67,1).NumberFormat = "@"
$ws.Cells.Item(67,1).Value = '2026-01-30'
$ws.Cells.Item(67,2).Value = '13:02:22'
$ws.Cells.Item(67,3).Value = '13:00'
$ws.Cells.Item(67,4).Value = 'Bathroom'
$ws.Cells.Item(67,5).Value = 'No Motion'
$ws.Cells.Item(67,6).Value = 'Inactive'
$ws.Cells.Item(68,1).NumberFormat = "@"
$ws.Cells.Item(68,1).Value = '2026-01-30'
$ws.Cells.Item(68,2).Value = '13:02:27'
$ws.Cells.Item(68,3).Value = '13:00'
$ws.Cells.Item(68,4).Value = 'Bathroom'
$ws.Cells.Item(68,5).Value = 'No Motion'
$ws.Cells.Item(68,6).Value = 'Inactive'
$ws.Cells.Item(69,1).NumberFormat = "@"
$ws.Cells.Item(69,1).Value = '2026-01-30'
$ws.Cells.Item(69,2).Value = '13:02:32'
$ws.Cells.Item(69,3).Value = '13:00'
$ws.Cells.Item(69,4).Value = 'Bathroom'
$ws.Cells.Item(69,5).Value = 'No Motion'
$ws.Cells.Item(69,6).Value = 'Inactive'
$ws.Cells.Item(70,1).NumberFormat = "@"
$ws.Cells.Item(70,1).Value = '2026-01-30'
$ws.Cells.Item(70,2).Value = '13:02:37'
$ws.Cells.Item(70,3).Value = '13:00'
$ws.Cells.Item(70,4).Value = 'Bathroom'
$ws.Cells.Item(70,5).Value = 'No Motion'
$ws.Cells.Item(70,6).Value = 'Inactive'
$ws.Cells.Item(71,1).NumberFormat = "@"
$ws.Cells.Item(71,1).Value = '2026-01-30'
$ws.Cells.Item(71,2).Value = '13:02:43'
$ws.Cells.Item(71,3).Value = '13:00'
$ws.Cells.Item(71,4).Value = 'Bathroom'
$ws.Cells.Item(71,5).Value = 'No Motion'
$ws.Cells.Item(71,6).Value = 'Inactive'
$ws.Cells.Item(72,1).NumberFormat = "@"
$ws.Cells.Item(72,1).Value = '2026-01-30'
$ws.Cells.Item(72,2).Value = '13:02:47'
$ws.Cells.Item(72,3).Value = '13:00'
$ws.Cells.Item(72,4).Value = 'Bathroom'
$ws.Cells.Item(72,5).Value = 'No Motion'
$ws.Cells.Item(72,6).Value = 'Inactive'
$ws.Cells.Item(73,1).NumberFormat = "@"
$ws.Cells.Item(73,1).Value = '2026-01-30'
$ws.Cells.Item(73,2).Value = '13:02:53'
$ws.Cells.Item(73,3).Value = '13:00'
$ws.Cells.Item(73,4).Value = 'Bathroom'
$ws.Cells.Item(73,5).Value = 'No Motion'
$ws.Cells.Item(73,6).Value = 'Inactive'
$ws.Cells.Item(74,1).NumberFormat = "@"
$ws.Cells.Item(74,1).Value = '2026-01-30'
$ws.Cells.Item(74,2).Value = '13:02:58'
$ws.Cells.Item(74,3).Value = '13:00'
$ws.Cells.Item(74,4).Value = 'Bathroom'
$ws.Cells.Item(74,5).Value = 'No Motion'
$ws.Cells.Item(74,6).Value = 'Inactive'

# --- Humidity: append 3 rows starting at row 53 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(53,1).NumberFormat = "@"
$ws.Cells.Item(53,1).Value = '2026-01-30'
$ws.Cells.Item(53,2).Value = '12:56:35'
$ws.Cells.Item(53,3).Value = '12:00'
$ws.Cells.Item(53,4).Value = 'Bathroom'
$ws.Cells.Item(53,5).NumberFormat = "@"
$ws.Cells.Item(53,5).Value = '86.8%'
$ws.Cells.Item(53,6).Value = 'Active'
$ws.Cells.Item(54,1).NumberFormat = "@"
$ws.Cells.Item(54,1).Value = '2026-01-30'
$ws.Cells.Item(54,2).Value = '12:56:36'
$ws.Cells.Item(54,3).Value = '12:00'
$ws.Cells.Item(54,4).Value = 'Bathroom'
$ws.Cells.Item(54,5).NumberFormat = "@"
$ws.Cells.Item(54,5).Value = '87.7%'
$ws.Cells.Item(54,6).Value = 'Active'
$ws.Cells.Item(55,1).NumberFormat = "@"
$ws.Cells.Item(55,1).Value = '2026-01-30'
$ws.Cells.Item(55,2).Value = '12:56:39'
$ws.Cells.Item(55,3).Value = '12:00'
$ws.Cells.Item(55,4).Value = 'Bathroom'
$ws.Cells.Item(55,5).NumberFormat = "@"
$ws.Cells.Item(55,5).Value = '87.7%'
$ws.Cells.Item(55,6).Value = 'Active'

# --- Temperature: append 3 rows starting at row 53 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(53,1).NumberFormat = "@"
$ws.Cells.Item(53,1).Value = '2026-01-30'
$ws.Cells.Item(53,2).Value = '12:56:36'
$ws.Cells.Item(53,3).Value = '12:00'
$ws.Cells.Item(53,4).Value = 'Bathroom'
$ws.Cells.Item(53,5).Value = '22.6C'
$ws.Cells.Item(53,6).Value = 'Active'
$ws.Cells.Item(54,1).NumberFormat = "@"
$ws.Cells.Item(54,1).Value = '2026-01-30'
$ws.Cells.Item(54,2).Value = '12:56:36'
$ws.Cells.Item(54,3).Value = '12:00'
$ws.Cells.Item(54,4).Value = 'Bathroom'
$ws.Cells.Item(54,5).Value = '22.6C'
$ws.Cells.Item(54,6).Value = 'Active'
$ws.Cells.Item(55,1).NumberFormat = "@"
$ws.Cells.Item(55,1).Value = '2026-01-30'
$ws.Cells.Item(55,2).Value = '12:56:39'
$ws.Cells.Item(55,3).Value = '12:00'
$ws.Cells.Item(55,4).Value = 'Bathroom'
$ws.Cells.Item(55,5).Value = '22.6C'
$ws.Cells.Item(55,6).Value = 'Active'

# --- Proximity: append 8 rows starting at row 25 ---
$ws = $wb.Worksheets.Item("Proximity")
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = '2026-01-30'
$ws.Cells.Item(25,2).Value = '12:56:36'
$ws.Cells.Item(25,3).Value = '12:00'
$ws.Cells.Item(25,4).Value = 'Bathroom Door'
$ws.Cells.Item(25,5).Value = 'EXIT'
$ws.Cells.Item(25,6).Value = 'User EXITED Bathroom'
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = '2026-01-30'
$ws.Cells.Item(26,2).Value = '13:02:02'
$ws.Cells.Item(26,3).Value = '13:00'
$ws.Cells.Item(26,4).Value = 'Bathroom Door'
$ws.Cells.Item(26,5).Value = 'ENTER'
$ws.Cells.Item(26,6).Value = 'User ENTERED Bathroom'
$ws.Cells.Item(27,1).NumberFormat = "@"
$ws.Cells.Item(27,1).Value = '2026-01-30'
$ws.Cells.Item(27,2).Value = '13:02:07'
$ws.Cells.Item(27,3).Value = '13:00'
$ws.Cells.Item(27,4).Value = 'Bathroom Door'
$ws.Cells.Item(27,5).Value = 'EXIT'
$ws.Cells.Item(27,6).Value = 'User EXITED Bathroom'
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = '2026-01-30'
$ws.Cells.Item(28,2).Value = '13:02:14'
$ws.Cells.Item(28,3).Value = '13:00'
$ws.Cells.Item(28,4).Value = 'Bathroom Door'
$ws.Cells.Item(28,5).Value = 'ENTER'
$ws.Cells.Item(28,6).Value = 'User ENTERED Bathroom'
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = '2026-01-30'
$ws.Cells.Item(29,2).Value = '13:02:20'
$ws.Cells.Item(29,3).Value = '13:00'
$ws.Cells.Item(29,4).Value = 'Bathroom Door'
$ws.Cells.Item(29,5).Value = 'EXIT'
$ws.Cells.Item(29,6).Value = 'User EXITED Bathroom'
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = '2026-01-30'
$ws.Cells.Item(30,2).Value = '13:02:30'
$ws.Cells.Item(30,3).Value = '13:00'
$ws.Cells.Item(30,4).Value = 'Bathroom Door'
$ws.Cells.Item(30,5).Value = 'ENTER'
$ws.Cells.Item(30,6).Value = 'User ENTERED Bathroom'
$ws.Cells.Item(31,1).NumberFormat = "@"
$ws.Cells.Item(31,1).Value = '2026-01-30'
$ws.Cells.Item(31,2).Value = '13:02:35'
$ws.Cells.Item(31,3).Value = '13:00'
$ws.Cells.Item(31,4).Value = 'Bathroom Door'
$ws.Cells.Item(31,5).Value = 'EXIT'
$ws.Cells.Item(31,6).Value = 'User EXITED Bathroom'
$ws.Cells.Item(32,1).NumberFormat = "@"
$ws.Cells.Item(32,1).Value = '2026-01-30'
$ws.Cells.Item(32,2).Value = '13:02:53'
$ws.Cells.Item(32,3).Value = '13:00'
$ws.Cells.Item(32,4).Value = 'Bathroom Door'
$ws.Cells.Item(32,5).Value = 'ENTER'
$ws.Cells.Item(32,6).Value = 'User ENTERED Bathroom'

# --- mmWave: append 8 rows starting at row 23 ---
$ws = $wb.Worksheets.Item("mmWave")
$ws.Cells.Item(23,1).NumberFormat = "@"
$ws.Cells.Item(23,1).Value = '2026-01-30'
$ws.Cells.Item(23,2).Value = '13:01:59'
$ws.Cells.Item(23,3).Value = '13:00'
$ws.Cells.Item(23,4).Value = 'Living Room'
$ws.Cells.Item(23,5).Value = 'FALL_DETECTED'
$ws.Cells.Item(23,6).Value = 'EMERGENCY'
$ws.Cells.Item(24,1).NumberFormat = "@"
$ws.Cells.Item(24,1).Value = '2026-01-30'
$ws.Cells.Item(24,2).Value = '13:02:00'
$ws.Cells.Item(24,3).Value = '13:00'
$ws.Cells.Item(24,4).Value = 'Living Room'
$ws.Cells.Item(24,5).Value = 'FALL_DETECTED'
$ws.Cells.Item(24,6).Value = 'EMERGENCY'
$ws.Cells.Item(25,1).NumberFormat = "@"
$ws.Cells.Item(25,1).Value = '2026-01-30'
$ws.Cells.Item(25,2).Value = '13:02:05'
$ws.Cells.Item(25,3).Value = '13:00'
$ws.Cells.Item(25,4).Value = 'Living Room'
$ws.Cells.Item(25,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(25,6).Value = 'Active'
$ws.Cells.Item(26,1).NumberFormat = "@"
$ws.Cells.Item(26,1).Value = '2026-01-30'
$ws.Cells.Item(26,2).Value = '13:02:15'
$ws.Cells.Item(26,3).Value = '13:00'
$ws.Cells.Item(26,4).Value = 'Living Room'
$ws.Cells.Item(26,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(26,6).Value = 'Active'
$ws.Cells.Item(27,1).NumberFormat = "@"
$ws.Cells.Item(27,1).Value = '2026-01-30'
$ws.Cells.Item(27,2).Value = '13:02:26'
$ws.Cells.Item(27,3).Value = '13:00'
$ws.Cells.Item(27,4).Value = 'Living Room'
$ws.Cells.Item(27,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(27,6).Value = 'Active'
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = '2026-01-30'
$ws.Cells.Item(28,2).Value = '13:02:36'
$ws.Cells.Item(28,3).Value = '13:00'
$ws.Cells.Item(28,4).Value = 'Living Room'
$ws.Cells.Item(28,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(28,6).Value = 'Active'
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = '2026-01-30'
$ws.Cells.Item(29,2).Value = '13:02:47'
$ws.Cells.Item(29,3).Value = '13:00'
$ws.Cells.Item(29,4).Value = 'Living Room'
$ws.Cells.Item(29,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(29,6).Value = 'Active'
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = '2026-01-30'
$ws.Cells.Item(30,2).Value = '13:02:57'
$ws.Cells.Item(30,3).Value = '13:00'
$ws.Cells.Item(30,4).Value = 'Living Room'
$ws.Cells.Item(30,5).Value = 'PRESENCE_DETECTED'
$ws.Cells.Item(30,6).Value = 'Active'

